# Update countries & provincias Spain
# Applies the COVID-19 data refresh described by the commit diff:
#  - timestamp label updated
#  - Estados Unidos / Brasil / Chequia / Venezuela / Bahamas / Barbados /
#    San Martin (Parte Holandesa) totals refreshed
#  - Japon's refreshed case count overtakes Polonia's (unchanged) count,
#    so the two countries swap list position (row 35 <-> row 36)
#  - Nueva Caledonia's refreshed numbers overtake Belice's (unchanged)
#    numbers, so the two swap list position (row 192 <-> row 193)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header / last-updated timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 02:04"

# --- Estados Unidos (row 4) -----------------------------------------------
$ws.Range("B4").Value = 1347125
$ws.Range("C4").Value = 25340
$ws.Range("D4").Value = 237205
$ws.Range("E4").Value = 1029884
$ws.Range("F4").Value = 16816
$ws.Range("G4").Value = 1421
$ws.Range("H4").Value = 80036

# --- Brasil (row 11) -------------------------------------------------------
$ws.Range("B11").Value = 156061
$ws.Range("C11").Value = 10169
$ws.Range("E11").Value = 83720
$ws.Range("G11").Value = 664
$ws.Range("H11").Value = 10656

# --- Japon / Polonia swap (rows 35-36) -------------------------------------
# Japon's updated totals now exceed Polonia's (unchanged) totals, so Japon
# takes row 35 and Polonia drops to row 36.
$ws.Range("A35").Value = "Japon"
$ws.Range("B35").Value = 15663
$ws.Range("C35").Value = 88
$ws.Range("D35").Value = 5906
$ws.Range("E35").Value = 9150
$ws.Range("F35").Value = 287
$ws.Range("G35").Value = 17
$ws.Range("H35").Value = 607

$ws.Range("A36").Value = "Polonia"
$ws.Range("B36").Value = 15651
$ws.Range("C36").Value = 285
$ws.Range("D36").Value = 5437
$ws.Range("E36").Value = 9429
$ws.Range("F36").Value = 160
$ws.Range("G36").Value = 9
$ws.Range("H36").Value = 785

# --- Chequia (row 50) -------------------------------------------------------
$ws.Range("D50").Value = 4447
$ws.Range("E50").Value = 3372

# --- Venezuela (row 128) ----------------------------------------------------
$ws.Range("B128").Value = 402
$ws.Range("C128").Value = 14
$ws.Range("E128").Value = 202

# --- Bahamas (row 165) -------------------------------------------------------
$ws.Range("D165").Value = 37
$ws.Range("E165").Value = 44

# --- Barbados (row 167) ------------------------------------------------------
$ws.Range("B167").Value = 84
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 57
$ws.Range("E167").Value = 20

# --- San Martin (Parte Holandesa) (row 170) ----------------------------------
$ws.Range("D170").Value = 46
$ws.Range("E170").Value = 15
$ws.Range("G170").Value = 1
$ws.Range("H170").Value = 15

# --- Nueva Caledonia / Belice swap (rows 192-193) ----------------------------
# Nueva Caledonia's updated totals now exceed Belice's (unchanged) totals, so
# Nueva Caledonia takes row 192 and Belice drops to row 193.
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2
